$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below receive a purely numeric-looking string (e.g. "552.16",
# "1.00", "0.0225") that must stay literal TEXT (matching the source
# workbook, which stores every Price/Volume cell as inline text). Excel
# auto-coerces such strings to Number on a plain .Value assignment, so we
# mark the cells as Text first, then restore the default "Normal" style
# afterwards so no stray number-format styling is left behind. A
# multi-area Range(...) does not reliably apply NumberFormat/Style to
# every area in this host, so loop cell-by-cell instead.
$textCells = @('D5','D6','D10','D13','D18','D19','D20','D21','D22','D23','D24','D25','D30','D31','D33','D39','D40','D41','D43','D46','D47','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.052.97'
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').Value = '2.420.49'
$ws.Range('E3').Value = '  +3.33%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '552.16'
$ws.Range('E5').Value = '  +2.12%  '
$ws.Range('D6').Value = '137.77'
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('E9').Value = '  +3.83%  '
$ws.Range('D10').Value = '5.76'
$ws.Range('E10').Value = '  +4.76%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').Value = '24.89'
$ws.Range('E13').Value = '  +4.70%  '
$ws.Range('D14').Value = '2.849.13'
$ws.Range('E14').Value = '  +3.32%  '
$ws.Range('D15').Value = '59.961.57'
$ws.Range('E15').Value = '  +3.66%  '
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('D17').Value = '2.411.52'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '11.34'
$ws.Range('E18').Value = '  +6.35%  '
$ws.Range('D19').Value = '4.38'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').Value = '331.89'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '65.15'
$ws.Range('E23').Value = '  +3.39%  '
$ws.Range('D24').Value = '0.170'
$ws.Range('E24').Value = '  +3.80%  '
$ws.Range('D25').Value = '8.57'
$ws.Range('E25').Value = '  +3.39%  '
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = '0.0₃0782'
$ws.Range('E28').Value = '  +6.47%  '
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = '169.50'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').Value = '6.26'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').Value = '18.68'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +0.79%  '
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '39.43'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '0.417'
$ws.Range('E40').Value = '  +10.78%  '
$ws.Range('D41').Value = '312.99'
$ws.Range('E41').Value = '  +8.68%  '
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').Value = '139.20'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('D46').Value = '19.51'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('D47').Value = '0.410'
$ws.Range('E47').Value = '  +7.49%  '
$ws.Range('D48').Value = '0.575'
$ws.Range('E48').Value = '  +1.51%  '
$ws.Range('D49').Value = '0.0225'
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('D50').Value = '17.70'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').Value = '11.06'
$ws.Range('E51').Value = '  -0.18%  '

# Drop the temporary Text number-format again so styling matches the
# original (unstyled) cells exactly.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
